{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the content changes described by the diff:\n//  - New title, author name, and author email (built from 4 concatenated parts)\n//  - Rewritten \"history\" body paragraph (shortened to 3 sentences about\n//    government/politics instead of the original 9-sentence essay)\n//  - Rewritten \"Summary\" paragraph body\n//  - A new empty paragraph appended at the very end of the document\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Paragraph indices (0-based), matching the original document layout:\n//   0 -> Title\n//   1 -> Author name\n//   2 -> Author email\n//   3 -> (blank separator paragraph)\n//   4 -> Long \"history\" body paragraph\n//   5 -> \"Summary\" heading\n//   6 -> Summary body paragraph\nconst pTitle = paragraphs.items[0];\nconst pAuthor = paragraphs.items[1];\nconst pEmail = paragraphs.items[2];\nconst pBody = paragraphs.items[4];\nconst pSummary = paragraphs.items[6];\n\n// Title\npTitle.insertText(\n  \"A Glimpse into the Art of Government: An Exploration of Civics and Politics\",\n  Word.InsertLocation.replace\n);\n\n// Author name\npAuthor.insertText(\"Clara Bennett\", Word.InsertLocation.replace);\n\n// Author email (was split across runs as \"marcuswalton56@abromail\" + \".\" + \"net\";\n// becomes \"clara\" + \".\" + \"bennett88@institute\" + \".\" + \"edu\")\npEmail.insertText(\"clara.bennett88@institute.edu\", Word.InsertLocation.replace);\n\n// Main body paragraph: replace the whole paragraph's text with the new\n// (much shorter) passage about government and politics.\npBody.insertText(\n  \"The exploration of government and politics provides a lens through which we can examine the interplay of power dynamics, decision-making processes, and the quest for justice.\" +\n    \" It encourages us to think critically, to challenge assumptions, and to recognize the interconnections between our actions and their broader implications.\" +\n    \" As we navigate the complexities of governance and political engagement, we gain a deeper appreciation for the rights and responsibilities that come with being a citizen, and we embrace the opportunity to contribute to a better future for ourselves and for generations to come.\",\n  Word.InsertLocation.replace\n);\n\n// Summary paragraph: replace the whole paragraph's text with the new conclusion.\npSummary.insertText(\n  \"In conclusion, government and politics are fundamental pillars of human society, shaping the structures, processes, and relationships that define how we live together.\" +\n    \" The study of government and politics provides a critical lens through which we can examine the intricacies of governance, the interplay of power, and the quest for a just and equitable society.\" +\n    \" It equips us with the knowledge, skills, and values necessary to navigate the political landscape, to participate effectively in the decision-making process, and to work towards a better future for all.\",\n  Word.InsertLocation.replace\n);\n\n// A new, empty paragraph is appended at the very end of the document body.\nbody.insertParagraph(\"\", Word.InsertLocation.end);\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the content changes described by the diff:\n#  - New title, author name, and author email (built from 4 concatenated parts)\n#  - Rewritten \"history\" body paragraph (shortened to 3 sentences about\n#    government/politics instead of the original 9-sentence essay)\n#  - Rewritten \"Summary\" paragraph body\n#  - A new empty paragraph appended at the very end of the document\n\n$d = $word.ActiveDocument\n\n# Paragraph indices (1-based), matching the original document layout:\n#   1 -> Title\n#   2 -> Author name\n#   3 -> Author email\n#   4 -> (blank separator paragraph)\n#   5 -> Long \"history\" body paragraph\n#   6 -> \"Summary\" heading\n#   7 -> Summary body paragraph\n$pTitle = $d.Paragraphs.Item(1)\n$pAuthor = $d.Paragraphs.Item(2)\n$pEmail = $d.Paragraphs.Item(3)\n$pBody = $d.Paragraphs.Item(5)\n$pSummary = $d.Paragraphs.Item(7)\n\n# Title\n$r = $pTitle.Range\n$d.Range($r.Start, $r.End).Text = \"A Glimpse into the Art of Government: An Exploration of Civics and Politics\"\n\n# Author name\n$r = $pAuthor.Range\n$d.Range($r.Start, $r.End).Text = \"Clara Bennett\"\n\n# Author email (was split across runs as \"marcuswalton56@abromail\" + \".\" + \"net\";\n# becomes \"clara\" + \".\" + \"bennett88@institute\" + \".\" + \"edu\")\n$r = $pEmail.Range\n$d.Range($r.Start, $r.End).Text = \"clara.bennett88@institute.edu\"\n\n# Main body paragraph: replace the whole paragraph's text with the new\n# (much shorter) passage about government and politics.\n$r = $pBody.Range\n$newBody = \"The exploration of government and politics provides a lens through which we can examine the interplay of power dynamics, decision-making processes, and the quest for justice.\" + \" It encourages us to think critically, to challenge assumptions, and to recognize the interconnections between our actions and their broader implications.\" + \" As we navigate the complexities of governance and political engagement, we gain a deeper appreciation for the rights and responsibilities that come with being a citizen, and we embrace the opportunity to contribute to a better future for ourselves and for generations to come.\"\n$d.Range($r.Start, $r.End).Text = $newBody\n\n# Summary paragraph: replace the whole paragraph's text with the new conclusion.\n$r = $pSummary.Range\n$newSummary = \"In conclusion, government and politics are fundamental pillars of human society, shaping the structures, processes, and relationships that define how we live together.\" + \" The study of government and politics provides a critical lens through which we can examine the intricacies of governance, the interplay of power, and the quest for a just and equitable society.\" + \" It equips us with the knowledge, skills, and values necessary to navigate the political landscape, to participate effectively in the decision-making process, and to work towards a better future for all.\"\n$d.Range($r.Start, $r.End).Text = $newSummary\n\n# A new, empty paragraph is appended at the very end of the document body.\n$endRange = $d.Content\n$endRange.Collapse(0)\n$endRange.InsertParagraphAfter()\n"}
